$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.8829803438055706
$ws.Range("C2").Value = 0.1755377646059628
$ws.Range("D2").Value = 0.07869776133452433
$ws.Range("E2").Value = 0.09941137218497431
$ws.Range("G2").Value = 1.498862315415892
$ws.Range("H2").Value = 1.208302126339021
$ws.Range("M2").Value = 0.3588998942875961
$ws.Range("B3").Value = 0.7915654878975147
$ws.Range("C3").Value = 0.1526226228488952
$ws.Range("D3").Value = 0.07133251497639037
$ws.Range("E3").Value = 0.09315340541534312
$ws.Range("G3").Value = 1.414789356164249
$ws.Range("H3").Value = 1.174302172605053
$ws.Range("M3").Value = 0.3237949660798591
$ws.Range("B4").Value = 0.7358813128300028
$ws.Range("C4").Value = 0.1385599860710158
$ws.Range("D4").Value = 0.06685335012912219
$ws.Range("E4").Value = 0.08937498524816334
$ws.Range("G4").Value = 1.36403068976486
$ws.Range("H4").Value = 1.154066190260011
$ws.Range("M4").Value = 0.3024417656031773
$ws.Range("B5").Value = 0.7133000115429127
$ws.Range("C5").Value = 0.1328307656777952
$ws.Range("D5").Value = 0.06503869437567289
$ws.Range("E5").Value = 0.08785101908407
$ws.Range("G5").Value = 1.343558948799028
$ws.Range("H5").Value = 1.145979045323827
$ws.Range("M5").Value = 0.2937898753701162
$ws.Range("B6").Value = 0.7095570339513699
$ws.Range("C6").Value = 0.1318795063042728
$ws.Range("D6").Value = 0.06473800890587711
$ws.Range("E6").Value = 0.08759890882613774
$ws.Range("G6").Value = 1.34017236820489
$ws.Range("H6").Value = 1.14464574337967
$ws.Range("M6").Value = 0.2923562114433693
$ws.Range("B7").Value = 0.7355763281404393
$ws.Range("C7").Value = 0.1384827146149235
$ws.Range("D7").Value = 0.06682883422968189
$ws.Range("E7").Value = 0.08935436907341199
$ws.Range("G7").Value = 1.363753744311339
$ws.Range("H7").Value = 1.153956482048699
$ws.Range("M7").Value = 0.3023248831059178
$ws.Range("B8").Value = 0.851367259767585
$ws.Range("C8").Value = 0.1676348284937887
$ws.Range("D8").Value = 0.07614914768433323
$ws.Range("E8").Value = 0.09724014635884259
$ws.Range("G8").Value = 1.469692485230212
$ws.Range("H8").Value = 1.19644499576799
$ws.Range("M8").Value = 0.3467533501333051
$ws.Range("B9").Value = 1.082041344766537
$ws.Range("C9").Value = 0.2248830162622824
$ws.Range("D9").Value = 0.09477868155960323
$ws.Range("E9").Value = 0.113227310595498
$ws.Range("G9").Value = 1.68447511285234
$ws.Range("H9").Value = 1.284929683595749
$ws.Range("M9").Value = 0.4355219891104696
$ws.Range("B10").Value = 1.25384579344518
$ws.Range("C10").Value = 0.2670286224909262
$ws.Range("D10").Value = 0.1086963214556107
$ws.Range("E10").Value = 0.12531461444744
$ws.Range("G10").Value = 1.846850386239112
$ws.Range("H10").Value = 1.35321438721985
$ws.Range("M10").Value = 0.5018151136100073
$ws.Range("B11").Value = 1.3325392911687
$ws.Range("C11").Value = 0.2862285267657967
$ws.Range("D11").Value = 0.1150813877566179
$ws.Range("E11").Value = 0.1308926057209661
$ws.Range("G11").Value = 1.921775173409515
$ws.Range("H11").Value = 1.385017961515189
$ws.Range("M11").Value = 0.5322232771562909
$ws.Range("B12").Value = 1.362418047789959
$ws.Range("C12").Value = 0.2935036089716334
$ws.Range("D12").Value = 0.1175072522520395
$ws.Range("E12").Value = 0.1330166537376343
$ws.Range("G12").Value = 1.950304408370357
$ws.Range("H12").Value = 1.397169780535251
$ws.Range("M12").Value = 0.5437753764940823
$ws.Range("B13").Value = 1.355979570248735
$ws.Range("C13").Value = 0.2919365815139088
$ws.Range("D13").Value = 0.1169844403290341
$ws.Range("E13").Value = 0.132558672059254
$ws.Range("G13").Value = 1.944153079391072
$ws.Range("H13").Value = 1.394547811245957
$ws.Range("M13").Value = 0.5412857537998548
$ws.Range("B14").Value = 1.334995836151961
$ws.Range("C14").Value = 0.286826959043168
$ws.Range("D14").Value = 0.1152808038384023
$ws.Range("E14").Value = 0.1310671141626329
$ws.Range("G14").Value = 1.924119121428646
$ws.Range("H14").Value = 1.386015512968811
$ws.Range("M14").Value = 0.5331729237431801
$ws.Range("B15").Value = 1.322153065296391
$ws.Range("C15").Value = 0.2836977687889544
$ws.Range("D15").Value = 0.1142383233173234
$ws.Range("E15").Value = 0.13015503736532
$ws.Range("G15").Value = 1.91186830368008
$ws.Range("H15").Value = 1.380803423299369
$ws.Range("M15").Value = 0.5282084565285743
$ws.Range("B16").Value = 1.248713979267052
$ws.Range("C16").Value = 0.2657744640592341
$ws.Range("D16").Value = 0.1082801460398031
$ws.Range("E16").Value = 0.1249517086413903
$ws.Range("G16").Value = 1.841975577400177
$ws.Range("H16").Value = 1.35115102126133
$ws.Range("M16").Value = 0.4998330130654551
$ws.Range("B17").Value = 1.203800710128064
$ws.Range("C17").Value = 0.2547865445093578
$ws.Range("D17").Value = 0.1046389572626367
$ws.Range("E17").Value = 0.1217802326643067
$ws.Range("G17").Value = 1.799373190159372
$ws.Range("H17").Value = 1.3331513205313
$ws.Range("M17").Value = 0.4824906694683762
$ws.Range("B18").Value = 1.178018390624572
$ws.Range("C18").Value = 0.2484691550442903
$ws.Range("D18").Value = 0.1025496973529556
$ws.Range("E18").Value = 0.1199635365263845
$ws.Range("G18").Value = 1.774968873667575
$ws.Range("H18").Value = 1.322867948436283
$ws.Range("M18").Value = 0.4725393796628197
$ws.Range("B19").Value = 1.169297579023464
$ws.Range("C19").Value = 0.2463306215979628
$ws.Range("D19").Value = 0.1018431711086691
$ws.Range("E19").Value = 0.1193497034802675
$ws.Range("G19").Value = 1.766722935856706
$ws.Range("H19").Value = 1.319398064851583
$ws.Range("M19").Value = 0.4691740591061375
$ws.Range("B20").Value = 1.208576554313311
$ws.Range("C20").Value = 0.255955957484673
$ws.Range("D20").Value = 0.1050260431154726
$ws.Range("E20").Value = 0.1221170677235861
$ws.Range("G20").Value = 1.803897961920171
$ws.Range("H20").Value = 1.335060203687732
$ws.Range("M20").Value = 0.4843343446397057
$ws.Range("B21").Value = 1.341157102739487
$ws.Range("C21").Value = 0.288327652819703
$ws.Range("D21").Value = 0.1157809847618267
$ws.Range("E21").Value = 0.1315048981395606
$ws.Range("G21").Value = 1.929999287744295
$ws.Range("H21").Value = 1.38851869512041
$ws.Range("M21").Value = 0.5355548405871389
$ws.Range("B22").Value = 1.428268865918994
$ws.Range("C22").Value = 0.3095108797539865
$ws.Range("D22").Value = 0.1228565765635778
$ws.Range("E22").Value = 0.1377092454313598
$ws.Range("G22").Value = 2.013330122578623
$ws.Range("H22").Value = 1.424090255292299
$ws.Range("M22").Value = 0.5692476732887286
$ws.Range("B23").Value = 1.38173274428226
$ws.Range("C23").Value = 0.2982023972226386
$ws.Range("D23").Value = 0.119075856613108
$ws.Range("E23").Value = 0.1343914465026614
$ws.Range("G23").Value = 1.968769545688701
$ws.Range("H23").Value = 1.405046422656085
$ws.Range("M23").Value = 0.5512449299002071
$ws.Range("B24").Value = 1.206417274070532
$ws.Range("C24").Value = 0.2554272668146496
$ws.Range("D24").Value = 0.1048510288248252
$ws.Range("E24").Value = 0.1219647639770542
$ws.Range("G24").Value = 1.801852037713871
$ws.Range("H24").Value = 1.334196995628787
$ws.Range("M24").Value = 0.483500759825489
$ws.Range("B25").Value = 1.019238509941204
$ws.Range("C25").Value = 0.2093835272936246
$ws.Range("D25").Value = 0.08969956281470104
$ws.Range("E25").Value = 0.1088440304754528
$ws.Range("G25").Value = 1.625588099416092
$ws.Range("H25").Value = 1.260426311733397
$ws.Range("M25").Value = 0.4113242722898178
